$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "pop3" validation block mirroring the existing pop1 (rows 2-11) and
# pop2 (rows 13-22) blocks, but using the text1..text9 placeholder values
# and the new "*This field must be digits" error message in column L.
#
# The brand-new shared strings are written first, and in the exact order
# they should appear in xl/sharedStrings.xml: pop3, text1..text9, then the
# new error message last.

$ws.Cells.Item(24,1).Value  = "pop3"

$ws.Cells.Item(24,6).Value  = "text1"
$ws.Cells.Item(25,6).Value  = "text2"
$ws.Cells.Item(26,6).Value  = "text3"
$ws.Cells.Item(27,6).Value  = "text4"
$ws.Cells.Item(28,6).Value  = "text5"
$ws.Cells.Item(29,6).Value  = "text6"
$ws.Cells.Item(30,6).Value  = "text7"
$ws.Cells.Item(31,6).Value  = "text8"
$ws.Cells.Item(32,6).Value  = "text9"

$ws.Cells.Item(24,12).Value = "*This field must be digits"
$ws.Cells.Item(25,12).Value = "*This field must be digits"
$ws.Cells.Item(26,12).Value = "*This field must be digits"
$ws.Cells.Item(27,12).Value = "*This field must be digits"
$ws.Cells.Item(28,12).Value = "*This field must be digits"
$ws.Cells.Item(29,12).Value = "*This field must be digits"
$ws.Cells.Item(30,12).Value = "*This field must be digits"
$ws.Cells.Item(31,12).Value = "*This field must be digits"
$ws.Cells.Item(32,12).Value = "*This field must be digits"

# Remaining cells all reuse shared strings that already exist in the
# workbook (same values as the pop1 block in rows 2-11).
$ws.Cells.Item(24,2).Value  = "LIVEHTA Automation - Test_NonOncology_Automation_3"
$ws.Cells.Item(24,3).Value  = "LIVEHTA Automation - Test_NonOncology_Automation_3_radio_button"
$ws.Cells.Item(24,4).Value  = "Clinical-Interventional"
$ws.Cells.Item(24,5).Value  = "total_record_number"
$ws.Cells.Item(24,7).Value  = "Clinical"
$ws.Cells.Item(24,8).Value  = "Clinical_radio_button"

$ws.Cells.Item(25,1).Value  = "pop3"
$ws.Cells.Item(25,4).Value  = "Clinical-RWE"
$ws.Cells.Item(25,5).Value  = "total_excluded_number"
$ws.Cells.Item(25,7).Value  = "Clinical"
$ws.Cells.Item(25,8).Value  = "Clinical_radio_button"

$ws.Cells.Item(26,1).Value  = "pop3"
$ws.Cells.Item(26,5).Value  = "total_screenedTiAb_number"

$ws.Cells.Item(27,1).Value  = "pop3"
$ws.Cells.Item(27,5).Value  = "total_excluded_screenedTiAb_number"

$ws.Cells.Item(28,1).Value  = "pop3"
$ws.Cells.Item(28,5).Value  = "fulltext_review"

$ws.Cells.Item(29,1).Value  = "pop3"
$ws.Cells.Item(29,5).Value  = "excluded_fulltext_review"

$ws.Cells.Item(30,1).Value  = "pop3"
$ws.Cells.Item(30,5).Value  = "total_greyliterature_number"

$ws.Cells.Item(31,1).Value  = "pop3"
$ws.Cells.Item(31,5).Value  = "original_studies"

$ws.Cells.Item(32,1).Value  = "pop3"
$ws.Cells.Item(32,5).Value  = "records_number"

# Update the view so the newly added block is the visible/selected area,
# mirroring what the author would have seen after typing the new rows.
$ws.Range("F24:L32").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 5
